# manu v9, new pore calcs
# 1) Refresh the "datetimeFigureOut" date placeholder text (5/26/2021 -> 6/7/2021)
#    on the slide master and on every slide layout.
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tf = $shp.TextFrame
            if ($tf.HasText -eq -1) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "5/26/2021") {
                    $tr.Text = "6/7/2021"
                }
            }
        }
    }
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateField $layouts.Item($li).Shapes
}

# 2) Slide 6 ("schematics"): the "Cover Crops" label loses its trailing "s"
#    (re-typed as two runs: "Cover " + "Crop").
$s6 = $p.Slides.Item(6)
$grp = $s6.Shapes.Item(2)
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    if ($shp.HasTextFrame -eq -1) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Cover Crops") {
            $tail = $tr.Characters(7, 5)
            $tail.Text = "Crop"
        }
    }
}

# 3) Slide 6: the lone "4" callout loses its now-redundant trailing endParaRPr
#    (simulated by clearing and retyping the run).
for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    if ($shp.Name -eq "TextBox 47") {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "4") {
            $tr.Delete()
            [void]$tr.InsertAfter("4")
        }
    }
}
